$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: game labels across columns B-H
$ws.Range("B1").Value = "Game 1"
$ws.Range("C1").Value = "Game 2"
$ws.Range("D1").Value = "Game 3"
$ws.Range("E1").Value = "Game 4"
$ws.Range("F1").Value = "Game 5"
$ws.Range("G1").Value = "Game 6"
$ws.Range("H1").Value = "Game 7"

# Column A labels, rows 3-6 first, then row 2 last (matches shared-string build order)
$ws.Range("A3").Value = "Deaths"
$ws.Range("A4").Value = "Last Hits"
$ws.Range("A5").Value = "Denies"
$ws.Range("A6").Value = "Game Time"
$ws.Range("A2").Value = "Win/Loss"

# Column A width matches the target layout (closest value the host's
# character->pixel quantization can reach to the authored 11.5546875)
$ws.Range("A1").ColumnWidth = 10.666666666666666

# Put selection on B2 as in target sheetView
$ws.Range("B2").Select()
